# Update computed market/profit columns (H-N) across all Leve sheets,
# reflecting a refreshed data pull from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3387.6
$ws.Range("I32").Value = 2996.6
$ws.Range("J32").Value = 3778.6
$ws.Range("K32").Value = 2996.6
$ws.Range("L32").Value = 3778.6
$ws.Range("M32").Value = -2670.6
$ws.Range("N32").Value = -4430.6
$ws.Range("H33").Value = 302.18182
$ws.Range("I33").Value = 314
$ws.Range("J33").Value = 249
$ws.Range("K33").Value = 314
$ws.Range("L33").Value = 249
$ws.Range("M33").Value = -85
$ws.Range("N33").Value = -707
$ws.Range("H64").Value = 4484.75
$ws.Range("J64").Value = 4816.6665
$ws.Range("L64").Value = 4816.6665
$ws.Range("N64").Value = -5312.6665
$ws.Range("H67").Value = 4484.75
$ws.Range("J67").Value = 4816.6665
$ws.Range("L67").Value = 4816.6665
$ws.Range("N67").Value = -6532.6665
$ws.Range("H69").Value = 22076.924
$ws.Range("J69").Value = 22750
$ws.Range("L69").Value = 68250
$ws.Range("N69").Value = -69998
$ws.Range("H72").Value = 22076.924
$ws.Range("J72").Value = 22750
$ws.Range("L72").Value = 204750
$ws.Range("N72").Value = -213486
$ws.Range("H92").Value = 655.8125
$ws.Range("I92").Value = 558.3077
$ws.Range("K92").Value = 558.3077
$ws.Range("M92").Value = 689.6923
$ws.Range("H98").Value = 1933.3334
$ws.Range("I98").Value = 1933.3334
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1933.3334
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -435.3334
$ws.Range("N98").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").ClearContents()
$ws.Range("H107").Value = 322.16666
$ws.Range("I107").Value = 322.16666
$ws.Range("K107").Value = 322.16666
$ws.Range("M107").Value = 1597.83334
$ws.Range("H122").Value = 1933.3334
$ws.Range("I122").Value = 1933.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5800.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3350.0002
$ws.Range("N122").ClearContents()
$ws.Range("H137").Value = 1252.7142
$ws.Range("I137").Value = 1033.3334
$ws.Range("J137").Value = 1647.6
$ws.Range("K137").Value = 3100.0002
$ws.Range("L137").Value = 4942.799999999999
$ws.Range("M137").Value = -550.0001999999999
$ws.Range("N137").Value = -10042.8
$ws.Range("H141").Value = 3124.75
$ws.Range("I141").Value = 1999.5
$ws.Range("K141").Value = 5998.5
$ws.Range("M141").Value = -818.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 816.7692
$ws.Range("I97").Value = 679.8889
$ws.Range("J97").Value = 1124.75
$ws.Range("K97").Value = 679.8889
$ws.Range("L97").Value = 1124.75
$ws.Range("M97").Value = -183.8889
$ws.Range("N97").Value = -2116.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2750
$ws.Range("I105").Value = 2750
$ws.Range("K105").Value = 2750
$ws.Range("M105").Value = -1003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2315.2
$ws.Range("I31").Value = 1816.75
$ws.Range("J31").Value = 2549.7646
$ws.Range("K31").Value = 1816.75
$ws.Range("L31").Value = 2549.7646
$ws.Range("M31").Value = -1521.75
$ws.Range("N31").Value = -3139.7646
$ws.Range("H34").Value = 2315.2
$ws.Range("I34").Value = 1816.75
$ws.Range("J34").Value = 2549.7646
$ws.Range("K34").Value = 1816.75
$ws.Range("L34").Value = 2549.7646
$ws.Range("M34").Value = -1614.75
$ws.Range("N34").Value = -2953.7646
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1225
$ws.Range("I13").Value = 966.6667
$ws.Range("K13").Value = 2900.0001
$ws.Range("M13").Value = -2732.0001
$ws.Range("H68").Value = 3227.875
$ws.Range("J68").Value = 3343.0667
$ws.Range("L68").Value = 10029.2001
$ws.Range("N68").Value = -11651.2001
$ws.Range("H71").Value = 3227.875
$ws.Range("J71").Value = 3343.0667
$ws.Range("L71").Value = 30087.6003
$ws.Range("N71").Value = -38199.6003
$ws.Range("H80").Value = 2988.5
$ws.Range("J80").Value = 2988.5
$ws.Range("L80").Value = 8965.5
$ws.Range("N80").Value = -10837.5
$ws.Range("H81").Value = 1955
$ws.Range("J81").Value = 1955
$ws.Range("L81").Value = 5865
$ws.Range("N81").Value = -8111
$ws.Range("H83").Value = 2988.5
$ws.Range("J83").Value = 2988.5
$ws.Range("L83").Value = 26896.5
$ws.Range("N83").Value = -36256.5
$ws.Range("H84").Value = 1955
$ws.Range("J84").Value = 1955
$ws.Range("L84").Value = 17595
$ws.Range("N84").Value = -28827
$ws.Range("H103").Value = 4633
$ws.Range("I103").Value = 3999
$ws.Range("J103").Value = 4950
$ws.Range("K103").Value = 11997
$ws.Range("L103").Value = 14850
$ws.Range("M103").Value = -11118
$ws.Range("N103").Value = -16608

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 1000000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H24").Value = 15000
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15346

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7386.6
$ws.Range("I7").Value = 9749
$ws.Range("J7").Value = 6796
$ws.Range("K7").Value = 9749
$ws.Range("L7").Value = 6796
$ws.Range("M7").Value = -9637
$ws.Range("N7").Value = -7020
$ws.Range("H55").Value = 1346.3334
$ws.Range("I55").Value = 1174.5454
$ws.Range("K55").Value = 1174.5454
$ws.Range("M55").Value = -1001.5454
$ws.Range("H93").Value = 652.6667
$ws.Range("I93").Value = 667.4545000000001
$ws.Range("K93").Value = 667.4545000000001
$ws.Range("M93").Value = 580.5454999999999
$ws.Range("H126").Value = 7386.6
$ws.Range("I126").Value = 9749
$ws.Range("J126").Value = 6796
$ws.Range("K126").Value = 29247
$ws.Range("L126").Value = 20388
$ws.Range("M126").Value = -26777
$ws.Range("N126").Value = -25328

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4884.636
$ws.Range("I122").Value = 3438.4285
$ws.Range("K122").Value = 10315.2855
$ws.Range("M122").Value = -7865.2855
$ws.Range("H126").Value = 1180.5
$ws.Range("I126").Value = 1196.5385
$ws.Range("K126").Value = 3589.6155
$ws.Range("M126").Value = -1119.6155
